$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row above the current row 409 (Excel-style "insert row
# above"), which shifts the existing rows 409:480 down to 410:481 and
# extends the used range to A1:T481.
$ws.Rows.Item(409).Insert()

# Populate the newly inserted row 409 with the new weekly price entry.
$ws.Range("A409").Value = 4
$ws.Range("B409").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C409").Value = "Los Lagos"
$ws.Range("D409").Value = 45218
$ws.Range("E409").Value = 10
$ws.Range("F409").Value = "Fruta"
$ws.Range("G409").Value = 100108
$ws.Range("H409").Value = "Tropicales y subtropicales"
$ws.Range("I409").Value = 100108005
$ws.Range("J409").Value = "Piña"
$ws.Range("K409").Value = "Caramelo"
$ws.Range("L409").Value = "Segunda"
$ws.Range("M409").Value = 100
$ws.Range("N409").Value = 25000
$ws.Range("O409").Value = 25000
$ws.Range("P409").Value = 25000
$ws.Range("Q409").Value = "$/caja 14 unidades"
$ws.Range("R409").Value = "Ecuador"
$ws.Range("S409").Value = 1786
$ws.Range("T409").Value = 14
